$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected (no way to know the original plaintext password from
# its stored hash), so unprotect it first in order to write the updated cells.
$ws.Unprotect()

# Update the confidential disclosure text (A7) - date changed 2021-04-08 -> 2021-04-09
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."

# Update holdings weight / percent-change figures
$ws.Range("D2").Value = 0.8437086284354224
$ws.Range("E2").Value = 0.003086816720257124

$ws.Range("D3").Value = 0.1562913715645776
$ws.Range("E3").Value = -0.008516941307165404

$ws.Range("E4").Value = 0.001273249462847792
